$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Insert()

$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

for ($r = 4; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = 30
    $ws.Cells.Item($r, 1).Font.Bold = $true
}

$ws.Range("A20").Value = 30

[void]$ws.Range("A3:A19").Select()

